{"js": "// Fix a typo introduced in the \"Not found epydoc\" troubleshooting note:\n// the command text reads \"$ sudo apt install python-epydocb\" but the\n// package is actually called \"python-epydoc\" (trailing stray \"b\").\n// Find the unique occurrence of \"epydocb\" in the document body and\n// replace it with \"epydoc\".\nconst results = context.document.body.search(\"epydocb\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"epydoc\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Fix a typo introduced in the \"Not found epydoc\" troubleshooting note:\n# the command text reads \"$ sudo apt install python-epydocb\" but the\n# package is actually called \"python-epydoc\" (trailing stray \"b\").\n# Find the unique occurrence of \"epydocb\" in the document and replace\n# it with \"epydoc\".\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"epydocb\"\n$find.Replacement.Text = \"epydoc\"\n$find.Execute([ref]\"epydocb\", $false, $false, $false, $false, $false, $true, 1, $false, \"epydoc\", 2)\n"}
